$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4529
$ws.Range("K3").Value = 4644
$ws.Range("K4").Value = 940
$ws.Range("K5").Value = 338
$ws.Range("K6").Value = 5247
$ws.Range("K7").Value = 15698

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 118
$ws.Range("K7").Value = 461
$ws.Range("K8").Value = 1056
$ws.Range("K9").Value = 67
$ws.Range("K10").Value = 83
$ws.Range("K11").Value = 311
$ws.Range("K15").Value = 161
$ws.Range("K17").Value = 30
$ws.Range("K18").Value = 107
$ws.Range("K19").Value = 473
$ws.Range("K20").Value = 363
$ws.Range("K21").Value = 45
$ws.Range("K23").Value = 161
$ws.Range("K25").Value = 76
$ws.Range("K29").Value = 831
$ws.Range("K30").Value = 56
$ws.Range("K31").Value = 173
$ws.Range("K33").Value = 658
$ws.Range("K36").Value = 200
$ws.Range("K37").Value = 532
$ws.Range("K42").Value = 582
$ws.Range("K46").Value = 35
$ws.Range("K47").Value = 100
$ws.Range("K48").Value = 200
$ws.Range("K49").Value = 90
$ws.Range("K51").Value = 196
$ws.Range("K52").Value = 416
$ws.Range("K53").Value = 210
$ws.Range("K54").Value = 295
$ws.Range("K55").Value = 176
$ws.Range("K58").Value = 9
$ws.Range("K60").Value = 102
$ws.Range("K63").Value = 46
$ws.Range("K65").Value = 360
$ws.Range("K66").Value = 51
$ws.Range("K67").Value = 604
$ws.Range("K76").Value = 219
$ws.Range("K79").Value = 391
$ws.Range("K83").Value = 337
$ws.Range("K85").Value = 705
$ws.Range("K86").Value = 105
$ws.Range("K87").Value = 28
$ws.Range("K88").Value = 182
$ws.Range("K89").Value = 224
$ws.Range("K94").Value = 196
$ws.Range("K95").Value = 275
$ws.Range("K98").Value = 79
$ws.Range("K101").Value = 15698

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 145
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 461

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 81
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 246
$ws.Range("K3").Value = 235
$ws.Range("K5").Value = 22
$ws.Range("K6").Value = 162
$ws.Range("K7").Value = 705

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 109
$ws.Range("K3").Value = 111
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 416

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 294
$ws.Range("K3").Value = 319
$ws.Range("K5").Value = 29
$ws.Range("K6").Value = 355
$ws.Range("K7").Value = 1056

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K4").Value = 17
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 337

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 181
$ws.Range("K3").Value = 245
$ws.Range("K7").Value = 658

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 92
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 148
$ws.Range("K3").Value = 176
$ws.Range("K7").Value = 532

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 109
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 360

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K2").Value = 14
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 59
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 173
$ws.Range("K6").Value = 175
$ws.Range("K7").Value = 604

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 48
$ws.Range("K6").Value = 150
$ws.Range("K7").Value = 295

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 239
$ws.Range("K3").Value = 293
$ws.Range("K6").Value = 234
$ws.Range("K7").Value = 831

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 200

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 147
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 473

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 46
$ws.Range("K6").Value = 121
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 156
$ws.Range("K3").Value = 182
$ws.Range("K6").Value = 217
$ws.Range("K7").Value = 582

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 45
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 128
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 391

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 115
$ws.Range("K7").Value = 363

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 200

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 35
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 42
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 9
